# Update "Horarios actualizados Línea 141 - 623"
# New scrape timestamp used throughout the workbook
$newTime = "03:00:15"

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: LP1912
# Row 6 (215_ALUAR) has already departed and is dropped from the list;
# the remaining rows shift up and a new arrival (11_ETCHEVERRY) is
# appended, growing the table from 3 to 4 data rows.
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: $newTime"
$ws1.Range("A3").Value = "Total filas: 4"

$ws1.Range("A6").Value = $newTime
$ws1.Range("B6").Value = "03:48"
$ws1.Range("C6").Value = "14_ABASTO"
$ws1.Range("D6").Value = 48
$ws1.Range("E6").Value = "LP1912"

$ws1.Range("A7").Value = $newTime
$ws1.Range("B7").Value = "04:01"
$ws1.Range("C7").Value = "81_EL PELIGRO"
$ws1.Range("D7").Value = 61
$ws1.Range("E7").Value = "LP1912"

$ws1.Range("A8").Value = $newTime
$ws1.Range("B8").Value = "04:46"
$ws1.Range("C8").Value = "215A_EL PATO"
$ws1.Range("D8").Value = 106
$ws1.Range("E8").Value = "LP1912"

$ws1.Range("A9").Value = $newTime
$ws1.Range("B9").Value = "04:53"
$ws1.Range("C9").Value = "11_ETCHEVERRY"
$ws1.Range("D9").Value = 113
$ws1.Range("E9").Value = "LP1912"

# ---------------------------------------------------------------------
# Sheet 2: LP1912-215
# Only arrival it tracks (215/215A) refreshed with the new scrape.
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = "Última actualización: $newTime"

$ws2.Range("A6").Value = $newTime
$ws2.Range("B6").Value = "04:46"
$ws2.Range("C6").Value = "215A_EL PATO"
$ws2.Range("D6").Value = 106
$ws2.Range("E6").Value = "LP1912"

# ---------------------------------------------------------------------
# Sheet 3: 6203-6173
# No arrivals tracked; only the scrape timestamp refreshes.
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = "Última actualización: $newTime"
